$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Year 3 / C1 / ANATOMY session 2 ---
$ws.Range("G3").Value = "Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("H3").Value = "49/221"

# --- Row 12: Year 3 / C1 / HISTOLOGY session 1 ---
$ws.Range("G12").Value = "mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"

# --- Row 25: Year 3 / C2 / ANATOMY session 2 ---
$ws.Range("G25").Value = "Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("H25").Value = "22/246"

# --- Row 34: Year 3 / C2 / HISTOLOGY session 1 ---
$ws.Range("G34").Value = "mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"

# --- Percentage-looking cells (L10, S15, S16) ---
# Assigning a "NN.N%" string directly via .Value makes the engine coerce it
# into a numeric percentage (and mint a new percent number-format style),
# which would alter the cell's style index away from its original one.
# To keep these as plain text (matching the source data) with their
# original style untouched, stage the literal text in a scratch cell via a
# text formula, then transfer only the computed value (not formatting).
$scratch = $ws.Range("ZZ1")

$scratch.Formula = "=""22.2%"""
$scratch.Copy()
$ws.Range("L10").PasteSpecial(-4163)

$scratch.Formula = "=""25.6%"""
$scratch.Copy()
$ws.Range("S15").PasteSpecial(-4163)

$scratch.Formula = "=""18.9%"""
$scratch.Copy()
$ws.Range("S16").PasteSpecial(-4163)

$scratch.Clear()
